$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 44096
$ws.Range("C2").Value = 277266
$ws.Range("D2").Value = 8486
$ws.Range("E2").Value = 38467
$ws.Range("F2").Value = 2200
$ws.Range("G2").Value = 13.87
$ws.Range("H2").Value = 25.93

# Row 3
$ws.Range("C3").Value = 162214
$ws.Range("D3").Value = 5218

# Row 4
$ws.Range("B4").Value = 44096
$ws.Range("C4").Value = 83193
$ws.Range("D4").Value = 2070
$ws.Range("E4").Value = 3412
$ws.Range("G4").Value = 5.93
$ws.Range("H4").Value = 3.34
$ws.Range("K4").Value = 57492
$ws.Range("L4").Value = 2036

# Row 8
$ws.Range("B8").Value = 44096
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "236253"
$ws.Range("C8").NumberFormat = "General"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "19153"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = 37018
$ws.Range("F8").Value = 5343
$ws.Range("G8").Value = 28.21
$ws.Range("K8").Value = 131200

# Row 10
$ws.Range("B10").Value = 44096
$ws.Range("C10").Value = 27124
$ws.Range("D10").Value = 424
$ws.Range("E10").Value = 6775
$ws.Range("G10").Value = 26.62
$ws.Range("H10").Value = 36.32
$ws.Range("K10").Value = 25449
$ws.Range("L10").Value = 424

# Row 14
$ws.Range("B14").Value = 44096
$ws.Range("C14").Value = 185148
$ws.Range("D14").Value = 2261
$ws.Range("E14").Value = 34086
$ws.Range("F14").Value = 622
$ws.Range("H14").Value = 27.51

# Row 15
$ws.Range("B15").Value = 44096
$ws.Range("C15").Value = 5146
$ws.Range("E15").Value = 966
$ws.Range("G15").Value = 20.64
$ws.Range("K15").Value = 4680

# Row 17
$ws.Range("B17").Value = 44096
$ws.Range("C17").Value = 131988
$ws.Range("D17").Value = 2304
$ws.Range("E17").Value = 33000
$ws.Range("F17").Value = 906
$ws.Range("G17").Value = 37.46
$ws.Range("H17").Value = 40.79
$ws.Range("K17").Value = 88083
$ws.Range("L17").Value = 2221

# Row 20
$ws.Range("B20").Value = 44095
$ws.Range("C20").Value = 53814
$ws.Range("D20").Value = 3601
$ws.Range("E20").Value = 7333
$ws.Range("G20").Value = 13.63
$ws.Range("H20").Value = 18.38

# Row 21
$ws.Range("B21").Value = 44096
$ws.Range("C21").Value = 117828
$ws.Range("D21").Value = 6616
$ws.Range("E21").Value = 25758
$ws.Range("F21").Value = 2531
$ws.Range("G21").Value = 21.86
$ws.Range("H21").Value = 38.26

# Row 24
$ws.Range("O24").Value = 'An error occurred. ... ValueError("invalid literal for int() with base 10: ''10,700''")'

# Row 26
$ws.Range("B26").Value = 44096
$ws.Range("C26").Value = 4316
$ws.Range("E26").Value = 232
$ws.Range("G26").Value = 5.87
$ws.Range("K26").Value = 3954

# Row 28
$ws.Range("B28").Value = 44096
$ws.Range("C28").Value = 127969
$ws.Range("D28").Value = 9328
$ws.Range("E28").Value = 12136
$ws.Range("H28").Value = 8.14

# Row 29
$ws.Range("B29").Value = 44096
$ws.Range("C29").Value = 167515
$ws.Range("D29").Value = 3085
$ws.Range("E29").Value = 14783
$ws.Range("G29").Value = 14.2
$ws.Range("K29").Value = 104100

# Row 30
$ws.Range("B30").Value = 44096
$ws.Range("C30").Value = 39232
$ws.Range("D30").Value = 432
$ws.Range("E30").Value = 5916
$ws.Range("G30").Value = 23.19
$ws.Range("K30").Value = 25508

# Row 32
$ws.Range("B32").Value = 44096
$ws.Range("C32").Value = 115319
$ws.Range("D32").Value = 1864
$ws.Range("E32").Value = 11638
$ws.Range("F32").Value = 488
$ws.Range("H32").Value = 32.02
$ws.Range("K32").Value = 44196
$ws.Range("L32").Value = 1524

# Row 33
$ws.Range("B33").Value = 44096
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "65044"
$ws.Range("C33").NumberFormat = "General"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "443"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "1366"
$ws.Range("E33").NumberFormat = "General"
$ws.Range("H33").Value = 1.58

# Row 34
$ws.Range("B34").Value = 44095
$ws.Range("C34").Value = 94021
$ws.Range("D34").Value = 2846
$ws.Range("E34").Value = 37272
$ws.Range("F34").Value = 1322
$ws.Range("G34").Value = 39.64
$ws.Range("H34").Value = 46.45

# Row 36
$ws.Range("B36").Value = 44096
$ws.Range("C36").Value = 91422
$ws.Range("D36").Value = 1979
$ws.Range("E36").Value = 16382
$ws.Range("G36").Value = 17.92
$ws.Range("H36").Value = 9.55

# Row 40
$ws.Range("B40").Value = 44096
$ws.Range("C40").Value = 41785
$ws.Range("D40").Value = 461
$ws.Range("E40").Value = 2257
$ws.Range("G40").Value = 7.03
$ws.Range("H40").Value = 7.71
$ws.Range("K40").Value = 32113
$ws.Range("L40").Value = 428

# Row 41
$ws.Range("B41").Value = 44095
$ws.Range("C41").Value = 784324
$ws.Range("D41").Value = 15071
$ws.Range("E41").Value = 23073
$ws.Range("F41").Value = 1140
$ws.Range("K41").Value = 541492
$ws.Range("L41").Value = 14723

# Row 43
$ws.Range("B43").Value = 44096
$ws.Range("C43").Value = 1721
$ws.Range("G43").Value = 10.45
$ws.Range("K43").Value = 1665

# Row 44
$ws.Range("B44").Value = 44096
$ws.Range("C44").Value = 27790
$ws.Range("D44").Value = 854
$ws.Range("G44").Value = 1.78

# Row 45
$ws.Range("B45").Value = 44096
$ws.Range("C45").Value = 679776
$ws.Range("D45").Value = 13416
$ws.Range("E45").Value = 91438
$ws.Range("F45").Value = 2354
$ws.Range("G45").Value = 13.45
$ws.Range("H45").Value = 17.55

# Row 48
$ws.Range("B48").Value = 44096
$ws.Range("C48").Value = 66053
$ws.Range("D48").Value = 2025
$ws.Range("E48").Value = 2873
$ws.Range("F48").Value = 129
$ws.Range("G48").Value = 5.21
$ws.Range("H48").Value = 6.49
$ws.Range("K48").Value = 55181
$ws.Range("L48").Value = 1987

# Row 49
$ws.Range("B49").Value = 44096
$ws.Range("C49").Value = 45147
$ws.Range("D49").Value = 765
$ws.Range("E49").Value = 1670
$ws.Range("G49").Value = 4.48
$ws.Range("H49").Value = 4.24
$ws.Range("K49").Value = 37262
$ws.Range("L49").Value = 755

# Row 51
$ws.Range("B51").Value = 44096
$ws.Range("C51").Value = 147070
$ws.Range("D51").Value = 8023
$ws.Range("E51").Value = 19973
$ws.Range("G51").Value = 25.68
$ws.Range("K51").Value = 77780

# Row 53
$ws.Range("B53").Value = 44096
$ws.Range("C53").Value = 308221
$ws.Range("D53").Value = 6677
$ws.Range("E53").Value = 80047
$ws.Range("F53").Value = 2774
$ws.Range("G53").Value = 25.97
$ws.Range("H53").Value = 41.55

# Row 56
$ws.Range("B56").Value = 44096
$ws.Range("C56").Value = 112626
$ws.Range("D56").Value = 3295
$ws.Range("E56").Value = 11005
$ws.Range("G56").Value = 9.77
$ws.Range("H56").Value = 13.47

# Row 58
$ws.Range("B58").Value = 44096
$ws.Range("C58").Value = 62731
$ws.Range("D58").Value = 1119
$ws.Range("E58").Value = 4962
$ws.Range("F58").Value = 132
$ws.Range("G58").Value = 11.66
$ws.Range("H58").Value = 12.93
$ws.Range("K58").Value = 42562
$ws.Range("L58").Value = 1020

# Row 59
$ws.Range("B59").Value = 44095
$ws.Range("C59").Value = 262133
$ws.Range("D59").Value = 6401
$ws.Range("E59").Value = 7545
$ws.Range("F59").Value = 591
$ws.Range("G59").Value = 4.81
$ws.Range("H59").Value = 9.81
$ws.Range("K59").Value = 156893
$ws.Range("L59").Value = 6022

# Row 38 (Texas) - previously empty/error row now fully populated
$ws.Range("B38").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("B38").Value = 44096
$ws.Range("C38").Value = 52768
$ws.Range("D38").Value = 14994
$ws.Range("E38").Value = 8716
$ws.Range("F38").Value = 1704
$ws.Range("G38").Value = 0.17
$ws.Range("H38").Value = 0.11
$ws.Range("I38").Value = $true
$ws.Range("O38").Value = "Success!"

Write-Host "Edit applied successfully"
